$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.455.65'
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").Value = '4.044.72'
$ws.Range("E3").Value = '  +0.82%  '
$ws.Range("E4").Value = '  -0.10%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '543.88'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.58%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '152.18'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("D7").Value = '4.038.98'
$ws.Range("E7").Value = '  +0.83%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.698'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +1.08%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.754'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.63%  '
$ws.Range("E11").Value = '  +0.67%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '54.00'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +12.65%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.0000331'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +1.69%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '10.96'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +2.06%  '
$ws.Range("D15").Value = '4.689.32'
$ws.Range("E15").Value = '  +0.78%  '
$ws.Range("D16").Value = '4.043.16'
$ws.Range("E16").Value = '  +0.86%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '14.38'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +1.91%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '20.72'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.96%  '
$ws.Range("E19").Value = '  +1.39%  '
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("D21").Value = '72.425.61'
$ws.Range("E21").Value = '  +1.29%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '448.97'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +4.37%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '98.08'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.54%  '
$ws.Range("E24").Value = '  +0.41%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '4.29'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +2.77%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '14.65'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +1.23%  '
$ws.Range("E27").Value = '  +14.01%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '11.31'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +1.77%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '10.85'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.84%  '
$ws.Range("E30").Value = '  +2.15%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '37.28'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.32%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '7.95'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +14.48%  '
$ws.Range("E33").Value = '  +3.30%  '
$ws.Range("E34").Value = '  +1.66%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '49.21'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +17.01%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '681.00'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.51%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '66.90'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.26%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.452'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +6.38%  '
$ws.Range("D39").Value = '0.0₃0889'
$ws.Range("E39").Value = '  +7.91%  '
$ws.Range("E40").Value = '  -2.86%  '
$ws.Range("E41").Value = '  -3.95%  '
$ws.Range("E42").Value = '  -1.64%  '
$ws.Range("E43").Value = '  +17.64%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("E45").Value = '  +2.28%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.08%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.152'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +0.95%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.69'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +2.51%  '
$ws.Range("B49").Value = 'LidoDAOToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '3.55'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +6.97%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '3.12'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +3.82%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '3.31'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -1.68%  '
